# MVP and SPL added
#
# Maurizio (row 14) receives one more MVP credit, which bumps his
# PointsxG / Total.
#
# Cerro earns a new SPL Bonus point. That raises his Total to 130,
# tying him with Federico (Damiano) at rank 17. The leaderboard is kept
# sorted by Total (descending), so the three rows 18-20 (Federico
# (Damiano) / Robi (Stoppi) / Cerro) get reshuffled: Cerro now leads at
# row 18, Federico (Damiano) slides to row 19, and Robi (Stoppi) slides
# to row 20 (rank 19).
#
# Lower down, Matteo (Riccardo) and Riccardo Ricci (both on Total 105,
# rank 21) swap places between rows 22 and 23 as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Stage the text ("Win Ratio") cells that need to move between rows via
# a real Copy so they stay literal shared-string text instead of being
# re-interpreted as a percentage number when assigned through .Value.
# A scratch row far outside the used range (A1:P75) holds the snapshots
# while we shuffle things around, then gets cleared again at the end.
# ---------------------------------------------------------------------
$scratchRow = 200

$ws.Cells.Item(18, 5).Copy($ws.Cells.Item($scratchRow, 1))   # old E18 "55%" -> scratch1
$ws.Cells.Item(19, 5).Copy($ws.Cells.Item($scratchRow, 2))   # old E19 "64%" -> scratch2
$ws.Cells.Item(20, 5).Copy($ws.Cells.Item($scratchRow, 3))   # old E20 "90%" -> scratch3
$ws.Cells.Item(22, 5).Copy($ws.Cells.Item($scratchRow, 4))   # old E22 "71%" -> scratch4
$ws.Cells.Item(23, 5).Copy($ws.Cells.Item($scratchRow, 5))   # old E23 "50%" -> scratch5

# Cache the old numeric rows 18-20 and 22-23 (plain scalars, so later
# writes into those rows cannot disturb the cached copies).
$old18 = @{ C=$ws.Cells.Item(18,3).Value(); D=$ws.Cells.Item(18,4).Value(); F=$ws.Cells.Item(18,6).Value(); G=$ws.Cells.Item(18,7).Value(); H=$ws.Cells.Item(18,8).Value(); I=$ws.Cells.Item(18,9).Value(); J=$ws.Cells.Item(18,10).Value(); K=$ws.Cells.Item(18,11).Value(); L=$ws.Cells.Item(18,12).Value(); M=$ws.Cells.Item(18,13).Value(); N=$ws.Cells.Item(18,14).Value(); O=$ws.Cells.Item(18,15).Value(); P=$ws.Cells.Item(18,16).Value() }
$old19 = @{ C=$ws.Cells.Item(19,3).Value(); D=$ws.Cells.Item(19,4).Value(); F=$ws.Cells.Item(19,6).Value(); G=$ws.Cells.Item(19,7).Value(); H=$ws.Cells.Item(19,8).Value(); I=$ws.Cells.Item(19,9).Value(); J=$ws.Cells.Item(19,10).Value(); K=$ws.Cells.Item(19,11).Value(); L=$ws.Cells.Item(19,12).Value(); M=$ws.Cells.Item(19,13).Value(); N=$ws.Cells.Item(19,14).Value(); O=$ws.Cells.Item(19,15).Value(); P=$ws.Cells.Item(19,16).Value() }
$old20 = @{ C=$ws.Cells.Item(20,3).Value(); D=$ws.Cells.Item(20,4).Value(); F=$ws.Cells.Item(20,6).Value(); G=$ws.Cells.Item(20,7).Value(); H=$ws.Cells.Item(20,8).Value(); I=$ws.Cells.Item(20,9).Value(); J=$ws.Cells.Item(20,10).Value(); K=$ws.Cells.Item(20,11).Value(); L=$ws.Cells.Item(20,12).Value(); M=$ws.Cells.Item(20,13).Value(); N=$ws.Cells.Item(20,14).Value(); O=$ws.Cells.Item(20,15).Value(); P=$ws.Cells.Item(20,16).Value() }

$old22 = @{ C=$ws.Cells.Item(22,3).Value(); D=$ws.Cells.Item(22,4).Value(); F=$ws.Cells.Item(22,6).Value(); G=$ws.Cells.Item(22,7).Value(); H=$ws.Cells.Item(22,8).Value(); I=$ws.Cells.Item(22,9).Value(); J=$ws.Cells.Item(22,10).Value(); K=$ws.Cells.Item(22,11).Value(); L=$ws.Cells.Item(22,12).Value(); M=$ws.Cells.Item(22,13).Value(); N=$ws.Cells.Item(22,14).Value(); O=$ws.Cells.Item(22,15).Value(); P=$ws.Cells.Item(22,16).Value() }
$old23 = @{ C=$ws.Cells.Item(23,3).Value(); D=$ws.Cells.Item(23,4).Value(); F=$ws.Cells.Item(23,6).Value(); G=$ws.Cells.Item(23,7).Value(); H=$ws.Cells.Item(23,8).Value(); I=$ws.Cells.Item(23,9).Value(); J=$ws.Cells.Item(23,10).Value(); K=$ws.Cells.Item(23,11).Value(); L=$ws.Cells.Item(23,12).Value(); M=$ws.Cells.Item(23,13).Value(); N=$ws.Cells.Item(23,14).Value(); O=$ws.Cells.Item(23,15).Value(); P=$ws.Cells.Item(23,16).Value() }

# ---------------------------------------------------------------------
# Row 14: Maurizio - MVP +1, PointsxG and Total recomputed.
# ---------------------------------------------------------------------
$ws.Cells.Item(14, 10).Value = 2       # J14 MVP            1 -> 2
$ws.Cells.Item(14, 14).Value = 10.38   # N14 PointsxG   10.19 -> 10.38
$ws.Cells.Item(14, 15).Value = 166     # O14 Total         163 -> 166

# ---------------------------------------------------------------------
# Row 18 becomes Cerro: old row 20's stats, plus the new SPL Bonus and
# the recomputed PointsxG / Total. Rank (B) stays 17.
# ---------------------------------------------------------------------
$ws.Cells.Item(18, 1).Value = "Cerro"
$ws.Cells.Item(18, 2).Value = 17
$ws.Cells.Item(18, 3).Value = $old20.C
$ws.Cells.Item(18, 4).Value = $old20.D
$ws.Cells.Item($scratchRow, 3).Copy($ws.Cells.Item(18, 5))    # E18 <- old E20 "90%"
$ws.Cells.Item(18, 6).Value = $old20.F
$ws.Cells.Item(18, 7).Value = $old20.G
$ws.Cells.Item(18, 8).Value = $old20.H
$ws.Cells.Item(18, 9).Value = $old20.I
$ws.Cells.Item(18, 10).Value = $old20.J
$ws.Cells.Item(18, 11).Value = 1        # K18 SPL Bonus   0 -> 1 (new)
$ws.Cells.Item(18, 12).Value = $old20.L
$ws.Cells.Item(18, 13).Value = $old20.M
$ws.Cells.Item(18, 14).Value = 13       # N18 PointsxG 12.7 -> 13
$ws.Cells.Item(18, 15).Value = 130      # O18 Total     127 -> 130
$ws.Cells.Item(18, 16).Value = $old20.P

# ---------------------------------------------------------------------
# Row 19 becomes Federico (Damiano): old row 18's stats, unchanged.
# ---------------------------------------------------------------------
$ws.Cells.Item(19, 1).Value = "Federico (Damiano)"
$ws.Cells.Item(19, 2).Value = 17
$ws.Cells.Item(19, 3).Value = $old18.C
$ws.Cells.Item(19, 4).Value = $old18.D
$ws.Cells.Item($scratchRow, 1).Copy($ws.Cells.Item(19, 5))    # E19 <- old E18 "55%"
$ws.Cells.Item(19, 6).Value = $old18.F
$ws.Cells.Item(19, 7).Value = $old18.G
$ws.Cells.Item(19, 8).Value = $old18.H
$ws.Cells.Item(19, 9).Value = $old18.I
$ws.Cells.Item(19, 10).Value = $old18.J
$ws.Cells.Item(19, 11).Value = $old18.K
$ws.Cells.Item(19, 12).Value = $old18.L
$ws.Cells.Item(19, 13).Value = $old18.M
$ws.Cells.Item(19, 14).Value = $old18.N
$ws.Cells.Item(19, 15).Value = $old18.O
$ws.Cells.Item(19, 16).Value = $old18.P

# ---------------------------------------------------------------------
# Row 20 becomes Robi (Stoppi): old row 19's stats, just rank (B) moves
# from 18 to 19 since Cerro now occupies rank 17 alongside Federico.
# ---------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value = "Robi (Stoppi)"
$ws.Cells.Item(20, 2).Value = 19
$ws.Cells.Item(20, 3).Value = $old19.C
$ws.Cells.Item(20, 4).Value = $old19.D
$ws.Cells.Item($scratchRow, 2).Copy($ws.Cells.Item(20, 5))    # E20 <- old E19 "64%"
$ws.Cells.Item(20, 6).Value = $old19.F
$ws.Cells.Item(20, 7).Value = $old19.G
$ws.Cells.Item(20, 8).Value = $old19.H
$ws.Cells.Item(20, 9).Value = $old19.I
$ws.Cells.Item(20, 10).Value = $old19.J
$ws.Cells.Item(20, 11).Value = $old19.K
$ws.Cells.Item(20, 12).Value = $old19.L
$ws.Cells.Item(20, 13).Value = $old19.M
$ws.Cells.Item(20, 14).Value = $old19.N
$ws.Cells.Item(20, 15).Value = $old19.O
$ws.Cells.Item(20, 16).Value = $old19.P

# Row 21 (Francesco) is untouched.

# ---------------------------------------------------------------------
# Row 22 becomes Riccardo Ricci: old row 23's stats, unchanged.
# ---------------------------------------------------------------------
$ws.Cells.Item(22, 1).Value = "Riccardo Ricci"
$ws.Cells.Item(22, 2).Value = 21
$ws.Cells.Item(22, 3).Value = $old23.C
$ws.Cells.Item(22, 4).Value = $old23.D
$ws.Cells.Item($scratchRow, 5).Copy($ws.Cells.Item(22, 5))    # E22 <- old E23 "50%"
$ws.Cells.Item(22, 6).Value = $old23.F
$ws.Cells.Item(22, 7).Value = $old23.G
$ws.Cells.Item(22, 8).Value = $old23.H
$ws.Cells.Item(22, 9).Value = $old23.I
$ws.Cells.Item(22, 10).Value = $old23.J
$ws.Cells.Item(22, 11).Value = $old23.K
$ws.Cells.Item(22, 12).Value = $old23.L
$ws.Cells.Item(22, 13).Value = $old23.M
$ws.Cells.Item(22, 14).Value = $old23.N
$ws.Cells.Item(22, 15).Value = $old23.O
$ws.Cells.Item(22, 16).Value = $old23.P

# ---------------------------------------------------------------------
# Row 23 becomes Matteo (Riccardo): old row 22's stats, unchanged.
# ---------------------------------------------------------------------
$ws.Cells.Item(23, 1).Value = "Matteo (Riccardo)"
$ws.Cells.Item(23, 2).Value = 21
$ws.Cells.Item(23, 3).Value = $old22.C
$ws.Cells.Item(23, 4).Value = $old22.D
$ws.Cells.Item($scratchRow, 4).Copy($ws.Cells.Item(23, 5))    # E23 <- old E22 "71%"
$ws.Cells.Item(23, 6).Value = $old22.F
$ws.Cells.Item(23, 7).Value = $old22.G
$ws.Cells.Item(23, 8).Value = $old22.H
$ws.Cells.Item(23, 9).Value = $old22.I
$ws.Cells.Item(23, 10).Value = $old22.J
$ws.Cells.Item(23, 11).Value = $old22.K
$ws.Cells.Item(23, 12).Value = $old22.L
$ws.Cells.Item(23, 13).Value = $old22.M
$ws.Cells.Item(23, 14).Value = $old22.N
$ws.Cells.Item(23, 15).Value = $old22.O
$ws.Cells.Item(23, 16).Value = $old22.P

# ---------------------------------------------------------------------
# Drop the scratch staging area again.
# ---------------------------------------------------------------------
$ws.Range($ws.Cells.Item($scratchRow, 1), $ws.Cells.Item($scratchRow, 5)).ClearContents()
